$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 23.666666
$ws.Range("I38").Value = 23.666666
$ws.Range("K38").Value = 70.99999800000001
$ws.Range("M38").Value = 301.000002
# Row 40
$ws.Range("H40").Value = 1819.8
$ws.Range("I40").Value = 1449.875
$ws.Range("J40").Value = 2242.5715
$ws.Range("K40").Value = 1449.875
$ws.Range("L40").Value = 2242.5715
$ws.Range("M40").Value = -1274.875
$ws.Range("N40").Value = -2592.5715
# Row 58
$ws.Range("H58").Value = 99.5
$ws.Range("I58").Value = 99.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 298.5
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -148.5
# Row 61
$ws.Range("H61").Value = 1078.3334
$ws.Range("I61").Value = 1078.3334
$ws.Range("K61").Value = 3235.0002
$ws.Range("M61").Value = -3063.0002
# Row 103
$ws.Range("H103").Value = 1066.6666
$ws.Range("I103").Value = 1150
$ws.Range("K103").Value = 3450
$ws.Range("M103").Value = -2864
# Row 107
$ws.Range("H107").Value = 45391.55
$ws.Range("I107").Value = 56638.5
$ws.Range("J107").Value = 403.75
$ws.Range("K107").Value = 56638.5
$ws.Range("L107").Value = 403.75
$ws.Range("M107").Value = -54718.5
$ws.Range("N107").Value = -4243.75
# Row 110
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180
# Row 132
$ws.Range("H132").Value = 806.4286
$ws.Range("I132").Value = 806.4286
$ws.Range("K132").Value = 2419.2858
$ws.Range("M132").Value = 110.7142000000003

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 629.7143
$ws.Range("I2").Value = 651.3333
$ws.Range("K2").Value = 651.3333
$ws.Range("M2").Value = -538.3333
# Row 5
$ws.Range("H5").Value = 433.25
$ws.Range("I5").Value = 433.25
$ws.Range("K5").Value = 433.25
$ws.Range("M5").Value = -321.25
# Row 32
$ws.Range("H32").Value = 6054.5557
$ws.Range("I32").Value = 6054.5557
$ws.Range("K32").Value = 6054.5557
$ws.Range("M32").Value = -5767.5557
# Row 41
$ws.Range("H41").Value = 1517.5
$ws.Range("I41").Value = 621
$ws.Range("K41").Value = 621
$ws.Range("M41").Value = -207
# Row 50
$ws.Range("H50").Value = 57155610
$ws.Range("I50").Value = 80000250
$ws.Range("J50").Value = 44000
$ws.Range("K50").Value = 80000250
$ws.Range("L50").Value = 44000
$ws.Range("M50").Value = -79999536
$ws.Range("N50").Value = -45428
# Row 74
$ws.Range("H74").Value = 1068.1428
$ws.Range("I74").Value = 1068.1428
$ws.Range("K74").Value = 1068.1428
$ws.Range("M74").Value = -194.1428000000001
# Row 77
$ws.Range("H77").Value = 1068.1428
$ws.Range("I77").Value = 1068.1428
$ws.Range("K77").Value = 5340.714
$ws.Range("M77").Value = -972.7139999999999
# Row 110
$ws.Range("I110").Value = 708.75
$ws.Range("K110").Value = 708.75
$ws.Range("M110").Value = 1336.25
# Row 116
$ws.Range("H116").Value = 629.7143
$ws.Range("I116").Value = 651.3333
$ws.Range("K116").Value = 651.3333
$ws.Range("M116").Value = 1642.6667
# Row 132
$ws.Range("H132").Value = 1705
$ws.Range("I132").Value = 1329.6154
$ws.Range("J132").Value = 3331.6667
$ws.Range("K132").Value = 3988.8462
$ws.Range("L132").Value = 9995.000100000001
$ws.Range("M132").Value = -1458.8462
$ws.Range("N132").Value = -15055.0001
# Row 139
$ws.Range("H139").Value = 71905
$ws.Range("J139").Value = 71905
$ws.Range("L139").Value = 71905
$ws.Range("N139").Value = -82185

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 629.7143
$ws.Range("I3").Value = 651.3333
$ws.Range("K3").Value = 651.3333
$ws.Range("M3").Value = -537.3333
# Row 4
$ws.Range("H4").Value = 433.25
$ws.Range("I4").Value = 433.25
$ws.Range("K4").Value = 433.25
$ws.Range("M4").Value = -318.25
# Row 22
$ws.Range("H22").Value = 910.2222
$ws.Range("I22").Value = 884.8570999999999
$ws.Range("K22").Value = 884.8570999999999
$ws.Range("M22").Value = -711.8570999999999
# Row 25
$ws.Range("H25").Value = 1471.6666
$ws.Range("I25").Value = 1207.5
$ws.Range("J25").Value = 2000
$ws.Range("K25").Value = 1207.5
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = -972.5
$ws.Range("N25").Value = -2470
# Row 105
$ws.Range("H105").Value = 920.3333
$ws.Range("I105").Value = 885.625
$ws.Range("J105").Value = 1198
$ws.Range("K105").Value = 885.625
$ws.Range("L105").Value = 1198
$ws.Range("M105").Value = 861.375
$ws.Range("N105").Value = -4692

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 97.42856999999999
$ws.Range("I7").Value = 197
$ws.Range("J7").Value = 22.75
$ws.Range("K7").Value = 197
$ws.Range("L7").Value = 22.75
$ws.Range("M7").Value = -84
$ws.Range("N7").Value = -248.75
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
# Row 38
$ws.Range("H38").Value = 31500
$ws.Range("J38").Value = 38500
$ws.Range("L38").Value = 38500
$ws.Range("N38").Value = -39254
# Row 46
$ws.Range("H46").Value = 31500
$ws.Range("J46").Value = 38500
$ws.Range("L46").Value = 38500
$ws.Range("N46").Value = -38922
# Row 58
$ws.Range("H58").Value = 1749.2
$ws.Range("I58").Value = 1582.3334
$ws.Range("K58").Value = 1582.3334
$ws.Range("M58").Value = -1379.3334
# Row 107
$ws.Range("H107").Value = 712.1667
$ws.Range("J107").Value = 766.3333
$ws.Range("L107").Value = 766.3333
$ws.Range("N107").Value = -4606.3333
# Row 122
$ws.Range("H122").Value = 3902.3845
$ws.Range("I122").Value = 1254.5714
$ws.Range("J122").Value = 6991.5
$ws.Range("K122").Value = 3763.7142
$ws.Range("L122").Value = 20974.5
$ws.Range("M122").Value = -1313.7142
$ws.Range("N122").Value = -25874.5
# Row 134
$ws.Range("H134").Value = 1413.1428
$ws.Range("J134").Value = 1147
$ws.Range("L134").Value = 3441
$ws.Range("N134").Value = -8511
# Row 136
$ws.Range("H136").Value = 1749.2
$ws.Range("I136").Value = 1582.3334
$ws.Range("K136").Value = 4747.0002
$ws.Range("M136").Value = -2197.0002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1339.7931
$ws.Range("I4").Value = 1348.25
$ws.Range("K4").Value = 4044.75
$ws.Range("M4").Value = -3932.75
# Row 99
$ws.Range("H99").Value = 3574.75
$ws.Range("I99").Value = 3574.75
$ws.Range("K99").Value = 10724.25
$ws.Range("M99").Value = -8478.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 97.12
$ws.Range("J2").Value = 102.5
$ws.Range("L2").Value = 102.5
$ws.Range("N2").Value = -328.5
# Row 9
$ws.Range("H9").Value = 293.5
$ws.Range("I9").Value = 293.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 293.5
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -123.5
# Row 17
$ws.Range("H17").Value = 836.3333
$ws.Range("J17").Value = 836.3333
$ws.Range("L17").Value = 836.3333
$ws.Range("N17").Value = -1172.3333
# Row 97
$ws.Range("H97").Value = 876.2857
$ws.Range("I97").Value = 372.33334
$ws.Range("K97").Value = 372.33334
$ws.Range("M97").Value = 123.66666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 42
$ws.Range("H42").Value = 23344.334
$ws.Range("I42").Value = 10025
$ws.Range("K42").Value = 10025
$ws.Range("M42").Value = -9462
# Row 49
$ws.Range("H49").Value = 23344.334
$ws.Range("I49").Value = 10025
$ws.Range("K49").Value = 10025
$ws.Range("M49").Value = -9878
# Row 55
$ws.Range("H55").Value = 585.5333000000001
$ws.Range("I55").Value = 443.8
$ws.Range("K55").Value = 443.8
$ws.Range("M55").Value = -270.8
# Row 93
$ws.Range("H93").Value = 47621864
$ws.Range("I93").Value = 66669500
$ws.Range("K93").Value = 66669500
$ws.Range("M93").Value = -66668252
# Row 136
$ws.Range("H136").Value = 4666.3335
$ws.Range("I136").Value = 3999.5
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 11998.5
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -9448.5
$ws.Range("N136").Value = -23100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1899.75
$ws.Range("I96").Value = 1539.8
$ws.Range("J96").Value = 2499.6667
$ws.Range("K96").Value = 1539.8
$ws.Range("L96").Value = 2499.6667
$ws.Range("M96").Value = -166.8
$ws.Range("N96").Value = -5245.6667
# Row 136
$ws.Range("H136").Value = 1720.3334
$ws.Range("I136").Value = 1795.4445
$ws.Range("K136").Value = 5386.333500000001
$ws.Range("M136").Value = -2836.333500000001
